$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook used to keep "property"/"state function"/"data node" tables as
# separate ad-hoc concepts. This commit unifies that naming -- the sheet that
# used to be called "Property1" is now "DataNode".
$ws.Name = "DataNode"

# Re-point the saved cursor/selection at D26 (beyond the used range), matching
# where the author last clicked before saving.
$ws.Range("D26").Select()
